# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between Week (A) and ASIN (old B)
#  - populate it with the Monday date for each forecast week
#  - shorten the week labels in column A from "W01".."W16" to "W1".."W16"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the current column B (ASIN), shifting
# ASIN..is_holiday_week one column to the right (B:I -> C:J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Make column B text-formatted so the ISO date strings aren't
# auto-converted into date serial numbers.
$ws.Columns.Item(2).NumberFormat = "@"

# Week start dates (Mondays), one per forecast week row.
$weekDates = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

foreach ($row in $weekDates.Keys) {
    $ws.Range("B$row").Value = $weekDates[$row]
}

# Shorten the week labels: "W01" -> "W1", "W02" -> "W2", ... "W16" stays "W16".
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Range("A$row")
    $week = $cell.Value()
    $cell.Value = $week -replace '^W0*(\d+)$', 'W$1'
}
